$d = $word.ActiveDocument

# The first occurrence below is wrapped in a bookmark; remove it up
# front (doing this after the text edit leaves the range stale and the
# bookmark un-resolvable), so no stray <w:bookmarkStart/> survives.
try {
    $bm = $d.Bookmarks.Item("_Hlk514861060")
    $bm.Delete()
} catch {
}

$newText = "Сазвежђе Orion: 16. до 25. јануара, 14. и 23. фебруара, 14. до 24. марта"
$marker = "Сазвежђе"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    $txt = $full.Text

    if ($txt.Contains($marker)) {
        # Exclude the trailing paragraph mark from the range so the
        # paragraph itself (and its pPr) survives the delete.
        $contentRange = $d.Range($full.Start, $full.End - 1)
        $contentRange.Delete()

        $p2 = $d.Paragraphs.Item($i)
        $insertPoint = $d.Range($p2.Range.Start, $p2.Range.Start)
        $insertPoint.Text = $newText
    }
}
